# Edit "Final Presentation Panel" - Slide 2 ("프로젝트 추진 개요")
#
# Re-positions a cluster of icon pictures / caption textboxes that sit in
# the left-hand "배경" (background) column of the slide (nudging them left
# and slightly down), and splits the long caption of the last bullet
# ("스마트 디바이스 및 프로그램과 연결을 통한 기능적 확장 요구") into two
# centred paragraphs.
#
# All target coordinates below are the EMU values taken straight from the
# target OOXML. The PowerPoint object model works in points
# (1 pt = 12700 EMU), so they are converted with ToPt(). A tiny +0.5 EMU
# bias is added before the division so that the single-precision
# (float32) storage used internally by Shape.Left/Top/Width/Height rounds
# back to the exact target EMU value instead of truncating one EMU short.

function ToPt($emu) {
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# --- "그림 62" (SD.png icon) -------------------------------------------------
$shp = Get-ShapeById $s 10
$shp.Left = ToPt(1380027)
$shp.Top  = ToPt(3531437)

# --- "TextBox 4" (시청각중복장애인 ... caption) ------------------------------
$shp = Get-ShapeById $s 5
$shp.Left   = ToPt(1972078)
$shp.Top    = ToPt(2942304)
$shp.Width  = ToPt(3095945)
$shp.Height = ToPt(502573)

# --- "TextBox 27" (값비싼 의사소통 보조기기 가격 caption) --------------------
$shp = Get-ShapeById $s 28
$shp.Left = ToPt(2410074)
$shp.Top  = ToPt(3734010)

# --- "TextBox 28" (보다 효율적인 의사소통 환경 구축 필요 caption) ------------
$shp = Get-ShapeById $s 29
$shp.Left = ToPt(2001952)
$shp.Top  = ToPt(4432416)

# --- "TextBox 29" (스마트 디바이스 ... caption) ------------------------------
$shp = Get-ShapeById $s 30
$shp.Left = ToPt(2264057)
$shp.Top  = ToPt(5277212)
# Split the single paragraph into two centred paragraphs.
$shp.TextFrame.TextRange.Text = "스마트 디바이스 및 프로그램과 `r연결을 통한 기능적 확장 요구"
# Re-assert the box height: spAutoFit recalculates it for the now
# two-line caption, but the target keeps the original box height.
$shp.Height = ToPt(502573)

# --- "그림 52" icon -----------------------------------------------------------
$shp = Get-ShapeById $s 53
$shp.Left = ToPt(1442561)
$shp.Top  = ToPt(5243456)

# --- "그림 53" icon -----------------------------------------------------------
$shp = Get-ShapeById $s 54
$shp.Left = ToPt(1427599)
$shp.Top  = ToPt(4370547)

# --- "그래픽 54" icon ---------------------------------------------------------
$shp = Get-ShapeById $s 55
$shp.Left = ToPt(1405196)
$shp.Top  = ToPt(2855730)
